$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf1"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.161357666666667
$ws.Range("H2").Value = 3.484073
$ws.Range("I2").Value = 0.1270850363824361
$ws.Range("J2").Value = 0.1270850363824361
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 87.94127933333334
$ws.Range("N2").Value = 263.823838
$ws.Range("O2").Value = 0.4109331243514438
$ws.Range("P2").Value = 0.4109331243514437
$ws.Range("Q2").Value = 102.1312789702416
$ws.Range("R2").Value = 919.181510732174
$ws.Range("S2").Value = 0.05222345105895137
$ws.Range("T2").Value = 0.05222345105895136

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf1"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.161357666666667
$ws.Range("H3").Value = 3.484073
$ws.Range("I3").Value = 0.1270850363824361
$ws.Range("J3").Value = 0.1270850363824361
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 52.441971
$ws.Range("N3").Value = 157.325913
$ws.Range("O3").Value = 0.2450515065683088
$ws.Range("P3").Value = 0.2450515065683087
$ws.Range("Q3").Value = 60.90388507596101
$ws.Range("R3").Value = 548.1349656836491
$ws.Range("S3").Value = 0.0311423796278043
$ws.Range("T3").Value = 0.03114237962780429

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf1"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.161357666666667
$ws.Range("H4").Value = 3.484073
$ws.Range("I4").Value = 0.1270850363824361
$ws.Range("J4").Value = 0.1270850363824361
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 54.667459
$ws.Range("N4").Value = 164.002377
$ws.Range("O4").Value = 0.255450795093328
$ws.Range("P4").Value = 0.255450795093328
$ws.Range("Q4").Value = 63.48847262683567
$ws.Range("R4").Value = 571.396253641521
$ws.Range("S4").Value = 0.03246397358835782
$ws.Range("T4").Value = 0.03246397358835781

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf1"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.161357666666667
$ws.Range("H5").Value = 3.484073
$ws.Range("I5").Value = 0.1270850363824361
$ws.Range("J5").Value = 0.1270850363824361
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 18.95316166666667
$ws.Range("N5").Value = 56.85948500000001
$ws.Range("O5").Value = 0.08856457398691947
$ws.Range("P5").Value = 0.08856457398691944
$ws.Range("Q5").Value = 22.01139960915611
$ws.Range("R5").Value = 198.102596482405
$ws.Range("S5").Value = 0.01125523210732261
$ws.Range("T5").Value = 0.01125523210732261

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf1"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.168173666666667
$ws.Range("H6").Value = 15.504521
$ws.Range("I6").Value = 0.565542861868062
$ws.Range("J6").Value = 0.565542861868062
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 87.94127933333334
$ws.Range("N6").Value = 263.823838
$ws.Range("O6").Value = 0.4109331243514438
$ws.Range("P6").Value = 0.4109331243514437
$ws.Range("Q6").Value = 454.4958040635109
$ws.Range("R6").Value = 4090.462236571599
$ws.Range("S6").Value = 0.2324002951820997
$ws.Range("T6").Value = 0.2324002951820997

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf1"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.168173666666667
$ws.Range("H7").Value = 15.504521
$ws.Range("I7").Value = 0.565542861868062
$ws.Range("J7").Value = 0.565542861868062
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 52.441971
$ws.Range("N7").Value = 157.325913
$ws.Range("O7").Value = 0.2450515065683088
$ws.Range("P7").Value = 0.2450515065683087
$ws.Range("Q7").Value = 271.029213550297
$ws.Range("R7").Value = 2439.262921952673
$ws.Range("S7").Value = 0.1385871303297215
$ws.Range("T7").Value = 0.1385871303297215

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf1"
$ws.Range("C8").Value = "Nrp1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.168173666666667
$ws.Range("H8").Value = 15.504521
$ws.Range("I8").Value = 0.565542861868062
$ws.Range("J8").Value = 0.565542861868062
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 54.667459
$ws.Range("N8").Value = 164.002377
$ws.Range("O8").Value = 0.255450795093328
$ws.Range("P8").Value = 0.255450795093328
$ws.Range("Q8").Value = 282.5309220273797
$ws.Range("R8").Value = 2542.778298246417
$ws.Range("S8").Value = 0.1444683737235526
$ws.Range("T8").Value = 0.1444683737235526

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf1"
$ws.Range("C9").Value = "Nrp1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.168173666666667
$ws.Range("H9").Value = 15.504521
$ws.Range("I9").Value = 0.565542861868062
$ws.Range("J9").Value = 0.565542861868062
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 18.95316166666667
$ws.Range("N9").Value = 56.85948500000001
$ws.Range("O9").Value = 0.08856457398691947
$ws.Range("P9").Value = 0.08856457398691944
$ws.Range("Q9").Value = 97.95323102574281
$ws.Range("R9").Value = 881.5790792316851
$ws.Range("S9").Value = 0.05008706263268815
$ws.Range("T9").Value = 0.05008706263268814

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fgf1"
$ws.Range("C10").Value = "Nrp1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.808898333333333
$ws.Range("H10").Value = 8.426695
$ws.Range("I10").Value = 0.3073721017495019
$ws.Range("J10").Value = 0.3073721017495019
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 87.94127933333334
$ws.Range("N10").Value = 263.823838
$ws.Range("O10").Value = 0.4109331243514438
$ws.Range("P10").Value = 0.4109331243514437
$ws.Range("Q10").Value = 247.0181129506011
$ws.Range("R10").Value = 2223.16301655541
$ws.Range("S10").Value = 0.1263093781103927
$ws.Range("T10").Value = 0.1263093781103927

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Fgf1"
$ws.Range("C11").Value = "Nrp1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.808898333333333
$ws.Range("H11").Value = 8.426695
$ws.Range("I11").Value = 0.3073721017495019
$ws.Range("J11").Value = 0.3073721017495019
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 52.441971
$ws.Range("N11").Value = 157.325913
$ws.Range("O11").Value = 0.2450515065683088
$ws.Range("P11").Value = 0.2450515065683087
$ws.Range("Q11").Value = 147.304164938615
$ws.Range("R11").Value = 1325.737484447535
$ws.Range("S11").Value = 0.07532199661078294
$ws.Range("T11").Value = 0.07532199661078293

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Fgf1"
$ws.Range("C12").Value = "Nrp1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2.808898333333333
$ws.Range("H12").Value = 8.426695
$ws.Range("I12").Value = 0.3073721017495019
$ws.Range("J12").Value = 0.3073721017495019
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 54.667459
$ws.Range("N12").Value = 164.002377
$ws.Range("O12").Value = 0.255450795093328
$ws.Range("P12").Value = 0.255450795093328
$ws.Range("Q12").Value = 153.5553344726683
$ws.Range("R12").Value = 1381.998010254015
$ws.Range("S12").Value = 0.07851844778141759
$ws.Range("T12").Value = 0.07851844778141756

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Fgf1"
$ws.Range("C13").Value = "Nrp1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2.808898333333333
$ws.Range("H13").Value = 8.426695
$ws.Range("I13").Value = 0.3073721017495019
$ws.Range("J13").Value = 0.3073721017495019
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 18.95316166666667
$ws.Range("N13").Value = 56.85948500000001
$ws.Range("O13").Value = 0.08856457398691947
$ws.Range("P13").Value = 0.08856457398691944
$ws.Range("Q13").Value = 53.23750421689724
$ws.Range("R13").Value = 479.1375379520751
$ws.Range("S13").Value = 0.0272222792469087
$ws.Range("T13").Value = 0.02722227924690869
